$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'72.018.00"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.41%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'4.021.09"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.86%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'534.35"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.48%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'150.00"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.77%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'4.016.35"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.69%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.693"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.81%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.02%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.754"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.41%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -3.86%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'54.37"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +11.48%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000325"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.89%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'10.82"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.91%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'4.660.89"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.89%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'4.010.95"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.04%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'14.20"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.15%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'20.80"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.50%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.46%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.69%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'71.993.89"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.39%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'432.51"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.47%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'98.50"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -3.10%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'3.59"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.46%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'14.74"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.78%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'4.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.99%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'4.29"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +25.64%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'11.43"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.65%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'10.84"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -2.03%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'5.95"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.61%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'37.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -1.75%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +21.92%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +2.43%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'50.23"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +17.90%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'13.62"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.57%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'673.79"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.58%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'68.69"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.75%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.455"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +3.96%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E40').Value = "'  -5.13%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'3.40"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +7.33%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -3.12%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.12%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'11.04"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +15.98%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0494"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -2.97%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.999"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.00%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.150"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -3.70%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.76%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'ApeXProtocol"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'3.42"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.47%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'Stacks"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'3.11"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.14%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'Maker"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.847.52"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +8.66%  "
$ws.Range('E51').Style = 'Normal'
